$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("constants")
$ws.Activate()

$ws.Range("E4").Value = 490
$ws.Range("G4").Value = 4180
$ws.Range("H4").Value = 405
$ws.Range("I4").Value = 323
$ws.Range("J4").Value = 280
$ws.Range("L4").Value = 2160
$ws.Range("M4").Value = 4080
$ws.Range("O4").Value = 200
$ws.Range("P4").Value = 200
$ws.Range("Q4").Value = 1000
$ws.Range("R4").Value = 728
$ws.Range("S4").Value = 914
$ws.Range("T4").Value = 662
$ws.Range("U4").Value = 499
$ws.Range("V4").Value = 432
$ws.Range("W4").Value = 3700
$ws.Range("X4").Value = 365
$ws.Range("Y4").Value = 520
$ws.Range("Z4").Value = 1710
$ws.Range("AA4").Value = 460

$ws.Range("R4").Select()
